$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 17
$ws.Range("H17").Value = 1378726.5
$ws.Range("J17").Value = 1421769.6
$ws.Range("L17").Value = 4265308.800000001
$ws.Range("N17").Value = -4265644.800000001

# Row 33
$ws.Range("H33").Value = 559.3333
$ws.Range("I33").Value = 291.63635
$ws.Range("J33").Value = 1295.5
$ws.Range("K33").Value = 291.63635
$ws.Range("L33").Value = 1295.5
$ws.Range("M33").Value = -62.63634999999999
$ws.Range("N33").Value = -1753.5

# Row 53
$ws.Range("H53").Value = 4077.0908
$ws.Range("I53").Value = 8212
$ws.Range("J53").Value = 631.3333
$ws.Range("K53").Value = 8212
$ws.Range("L53").Value = 631.3333
$ws.Range("M53").Value = -7575
$ws.Range("N53").Value = -1905.3333

# Row 86
$ws.Range("H86").Value = 24472.861
$ws.Range("I86").Value = 4779.615
$ws.Range("J86").Value = 40473.625
$ws.Range("K86").Value = 4779.615
$ws.Range("L86").Value = 40473.625
$ws.Range("M86").Value = -3656.615
$ws.Range("N86").Value = -42719.625

# Row 88
$ws.Range("H88").Value = 1275.3125
$ws.Range("I88").Value = 1065
$ws.Range("J88").Value = 1370.909
$ws.Range("K88").Value = 1065
$ws.Range("L88").Value = 1370.909
$ws.Range("M88").Value = -659
$ws.Range("N88").Value = -2182.909

# Row 89
$ws.Range("H89").Value = 24472.861
$ws.Range("I89").Value = 4779.615
$ws.Range("J89").Value = 40473.625
$ws.Range("K89").Value = 23898.075
$ws.Range("L89").Value = 202368.125
$ws.Range("M89").Value = -18282.075
$ws.Range("N89").Value = -213600.125

# Row 91
$ws.Range("H91").Value = 1275.3125
$ws.Range("I91").Value = 1065
$ws.Range("J91").Value = 1370.909
$ws.Range("K91").Value = 1065
$ws.Range("L91").Value = 1370.909
$ws.Range("M91").Value = 339
$ws.Range("N91").Value = -4178.909

# Row 137
$ws.Range("H137").Value = 7042.902
$ws.Range("I137").Value = 9351.971
$ws.Range("J137").Value = 2424.7646
$ws.Range("K137").Value = 28055.913
$ws.Range("L137").Value = 7274.293799999999
$ws.Range("M137").Value = -25505.913
$ws.Range("N137").Value = -12374.2938

$ws = $wb.Worksheets.Item("ARM")
# Row 132
$ws.Range("H132").Value = 4964
$ws.Range("I132").Value = 724.5
$ws.Range("K132").Value = 2173.5
$ws.Range("M132").Value = 356.5

$ws = $wb.Worksheets.Item("BSM")
# Row 86
$ws.Range("H86").Value = 5607.269
$ws.Range("I86").Value = 6379.5264
$ws.Range("J86").Value = 3511.1428
$ws.Range("K86").Value = 6379.5264
$ws.Range("L86").Value = 3511.1428
$ws.Range("M86").Value = -5256.5264
$ws.Range("N86").Value = -5757.1428

# Row 89
$ws.Range("H89").Value = 5607.269
$ws.Range("I89").Value = 6379.5264
$ws.Range("J89").Value = 3511.1428
$ws.Range("K89").Value = 31897.632
$ws.Range("L89").Value = 17555.714
$ws.Range("M89").Value = -26281.632
$ws.Range("N89").Value = -28787.714

$ws = $wb.Worksheets.Item("CRP")
# Row 31
$ws.Range("H31").Value = 11516.167
$ws.Range("J31").Value = 6999.5
$ws.Range("L31").Value = 6999.5
$ws.Range("N31").Value = -7589.5

# Row 34
$ws.Range("H34").Value = 11516.167
$ws.Range("J34").Value = 6999.5
$ws.Range("L34").Value = 6999.5
$ws.Range("N34").Value = -7403.5

# Row 58
$ws.Range("H58").Value = 2940.9524
$ws.Range("I58").Value = 1942.0769
$ws.Range("K58").Value = 1942.0769
$ws.Range("M58").Value = -1739.0769

# Row 120
$ws.Range("H120").Value = 63423
$ws.Range("I120").Value = 60296
$ws.Range("K120").Value = 60296
$ws.Range("M120").Value = -56667

# Row 132
$ws.Range("H132").Value = 2309.5454
$ws.Range("I132").Value = 2376.7778
$ws.Range("K132").Value = 7130.3334
$ws.Range("M132").Value = -4600.3334

# Row 134
$ws.Range("H134").Value = 6089.7144
$ws.Range("I134").Value = 4835.2856
$ws.Range("K134").Value = 14505.8568
$ws.Range("M134").Value = -11970.8568

# Row 136
$ws.Range("H136").Value = 2940.9524
$ws.Range("I136").Value = 1942.0769
$ws.Range("K136").Value = 5826.2307
$ws.Range("M136").Value = -3276.2307

# Row 141
$ws.Range("H141").Value = 101397
$ws.Range("J141").Value = 104071.836
$ws.Range("L141").Value = 104071.836
$ws.Range("N141").Value = -114431.836

$ws = $wb.Worksheets.Item("GSM")
# Row 80
$ws.Range("H80").Value = 14657.417
$ws.Range("I80").Value = 15808.182
$ws.Range("J80").Value = 1999
$ws.Range("K80").Value = 15808.182
$ws.Range("L80").Value = 1999
$ws.Range("M80").Value = -14810.182
$ws.Range("N80").Value = -3995

# Row 83
$ws.Range("H83").Value = 14657.417
$ws.Range("I83").Value = 15808.182
$ws.Range("J83").Value = 1999
$ws.Range("K83").Value = 79040.91
$ws.Range("L83").Value = 9995
$ws.Range("M83").Value = -74048.91
$ws.Range("N83").Value = -19979

# Row 132
$ws.Range("H132").Value = 3906.5
$ws.Range("I132").Value = 3906.5
$ws.Range("J132").Value = 0
$ws.Range("K132").Value = 11719.5
$ws.Range("L132").Value = 0
$ws.Range("M132").Value = -9189.5
$ws.Range("N132").ClearContents()

$ws = $wb.Worksheets.Item("LTW")
# Row 136
$ws.Range("H136").Value = 4713.814
$ws.Range("I136").Value = 3510
$ws.Range("J136").Value = 5580.56
$ws.Range("K136").Value = 10530
$ws.Range("L136").Value = 16741.68
$ws.Range("M136").Value = -7980
$ws.Range("N136").Value = -21841.68

$ws = $wb.Worksheets.Item("WVR")
# Row 10
$ws.Range("H10").Value = 6
$ws.Range("J10").Value = 6
$ws.Range("L10").Value = 6
$ws.Range("N10").Value = -344

# Row 62
$ws.Range("H62").Value = 300212
$ws.Range("I62").Value = 636000.7
$ws.Range("J62").Value = 12393.143
$ws.Range("K62").Value = 636000.7
$ws.Range("L62").Value = 12393.143
$ws.Range("M62").Value = -635376.7
$ws.Range("N62").Value = -13641.143

# Row 65
$ws.Range("H65").Value = 300212
$ws.Range("I65").Value = 636000.7
$ws.Range("J65").Value = 12393.143
$ws.Range("K65").Value = 3180003.5
$ws.Range("L65").Value = 61965.715
$ws.Range("M65").Value = -3176883.5
$ws.Range("N65").Value = -68205.715

# Row 96
$ws.Range("H96").Value = 16669462
$ws.Range("I96").Value = 33334090
$ws.Range("J96").Value = 4833.3335
$ws.Range("K96").Value = 33334090
$ws.Range("L96").Value = 4833.3335
$ws.Range("M96").Value = -33332717
$ws.Range("N96").Value = -7579.3335

# Row 132
$ws.Range("H132").Value = 14726.219
$ws.Range("I132").Value = 17742.727
$ws.Range("J132").Value = 8089.9
$ws.Range("K132").Value = 53228.181
$ws.Range("L132").Value = 24269.7
$ws.Range("M132").Value = -50698.181
$ws.Range("N132").Value = -29329.7

# Row 136
$ws.Range("H136").Value = 336277.25
$ws.Range("I136").Value = 386356.2
$ws.Range("J136").Value = 2417.6667
$ws.Range("K136").Value = 1159068.6
$ws.Range("L136").Value = 7253.000100000001
$ws.Range("M136").Value = -1156518.6
$ws.Range("N136").Value = -12353.0001
